# Apply cryptocurrency price/volume updates to D (Price) and E (Volume(1h)) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.145.73"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "2.284.79"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.08"
$ws.Range("E5").Value = "  -4.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.33"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "2.284.16"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("E10").Value = "  -5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.47"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").Value = "2.689.72"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "58.073.43"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "2.287.32"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.52"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  -5.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.55"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.38"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.63"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.99"
$ws.Range("E27").Value = "  -4.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.28"
$ws.Range("E28").Value = "  -6.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.52"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  -5.75%  "
$ws.Range("D31").Value = "0.0₃0719"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.75"
$ws.Range("E32").Value = "  -5.53%  "
$ws.Range("E33").Value = "  -6.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.380"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.73"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.89"
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.24"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.49"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.79"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.75"
$ws.Range("E43").Value = "  -9.79%  "
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.552"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.11"
$ws.Range("E48").Value = "  -6.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0210"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  -0.56%  "
